# "Changed a few things": add a new TA-assignment row (row 13) for
# Jacob Brown, mirroring the existing rows (e.g. row 6 / row 11) that
# already use him as the assigned TA, then move the selection onto the
# new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: # | TA assigned | TA Email
$ws.Range("A13").Value = 10
$ws.Range("B13").Value = "Jacob Brown"
$ws.Range("C13").Value = "jacbrow@okstate.edu"

# Make the e-mail a real clickable mailto: hyperlink, like the other
# "TA Email" cells in the sheet.
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:jacbrow@okstate.edu")

# Restore the same visual formatting the other hyperlinked e-mail cells
# use (Hyperlinks.Add() re-stamps its own format on the cell).
$ws.Range("C6").Copy()
$ws.Range("C13").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# The new row is what's now selected in the sheet.
$ws.Range("A13:E13").Select()
